$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 46

$ws.Cells.Item($row, 1).Value = "60OWJY"
$ws.Cells.Item($row, 2).Value = "Engranaje de cinta Epson"
$ws.Cells.Item($row, 3).Value = "TM U220"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 100000
$ws.Cells.Item($row, 6).Value = 1
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Formula = "=(E46-D46)*G46"
$ws.Cells.Item($row, 9).Formula = "=D46*F46"
$ws.Cells.Item($row, 10).Value = 0
